$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format for numeric-looking value columns (D, E, G)
# so Excel stores them as text, matching the original inlineStr string cells,
# instead of auto-converting to numbers/percentages/dates.
$textCols = @("D", "E", "G")
foreach ($col in $textCols) {
    $ws.Range($col + "2:" + $col + "51").NumberFormat = "@"
}

$ws.Range("D2").Value = "312.05"
$ws.Range("E2").Value = "1.06%"
$ws.Range("G2").Value = "10"
$ws.Range("D3").Value = "39.52"
$ws.Range("E3").Value = "2.68%"
$ws.Range("G3").Value = "10"
$ws.Range("D4").Value = "5.155"
$ws.Range("E4").Value = "1.04%"
$ws.Range("G4").Value = "10"
$ws.Range("D5").Value = "0.08187"
$ws.Range("E5").Value = "0.79%"
$ws.Range("G5").Value = "10"
$ws.Range("D6").Value = "1.983"
$ws.Range("E6").Value = "1.15%"
$ws.Range("G6").Value = "10"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "8.146"
$ws.Range("E7").Value = "2.62%"
$ws.Range("G7").Value = "10"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.226"
$ws.Range("E8").Value = "0.49%"
$ws.Range("G8").Value = "10"
$ws.Range("D9").Value = "0.9247"
$ws.Range("E9").Value = "-0.32%"
$ws.Range("G9").Value = "10"
$ws.Range("D10").Value = "0.1396"
$ws.Range("E10").Value = "-3.06%"
$ws.Range("G10").Value = "10"
$ws.Range("D11").Value = "0.1929"
$ws.Range("E11").Value = "-1.71%"
$ws.Range("G11").Value = "10"
$ws.Range("D12").Value = "0.09090"
$ws.Range("E12").Value = "-0.19%"
$ws.Range("G12").Value = "10"
$ws.Range("D13").Value = "0.03500"
$ws.Range("E13").Value = "-0.24%"
$ws.Range("G13").Value = "10"
$ws.Range("D14").Value = "0.09811"
$ws.Range("E14").Value = "-0.01%"
$ws.Range("G14").Value = "10"
$ws.Range("D15").Value = "0.001399"
$ws.Range("E15").Value = "-0.57%"
$ws.Range("G15").Value = "10"
$ws.Range("D16").Value = "0.006015"
$ws.Range("E16").Value = "-1.40%"
$ws.Range("G16").Value = "10"
$ws.Range("D17").Value = "3.673"
$ws.Range("E17").Value = "0.50%"
$ws.Range("G17").Value = "10"
$ws.Range("D18").Value = "3.353"
$ws.Range("E18").Value = "-3.83%"
$ws.Range("G18").Value = "10"
$ws.Range("D19").Value = "0.3463"
$ws.Range("E19").Value = "0.44%"
$ws.Range("G19").Value = "10"
$ws.Range("D20").Value = "0.1350"
$ws.Range("E20").Value = "1.16%"
$ws.Range("G20").Value = "10"
$ws.Range("D21").Value = "4.655"
$ws.Range("E21").Value = "-3.02%"
$ws.Range("G21").Value = "10"
$ws.Range("D22").Value = "0.2416"
$ws.Range("E22").Value = "-1.53%"
$ws.Range("G22").Value = "10"
$ws.Range("D23").Value = "0.04362"
$ws.Range("E23").Value = "-1.61%"
$ws.Range("G23").Value = "10"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").Value = "0.85%"
$ws.Range("G24").Value = "10"
$ws.Range("D25").Value = "0.004805"
$ws.Range("E25").Value = "-0.62%"
$ws.Range("G25").Value = "10"
$ws.Range("D26").Value = "0.0001295"
$ws.Range("E26").Value = "-0.38%"
$ws.Range("G26").Value = "10"
$ws.Range("D27").Value = "0.0003986"
$ws.Range("E27").Value = "-10.37%"
$ws.Range("G27").Value = "10"
$ws.Range("G28").Value = "10"
$ws.Range("G29").Value = "10"
$ws.Range("G30").Value = "10"
$ws.Range("G31").Value = "10"
$ws.Range("G32").Value = "10"
$ws.Range("G33").Value = "10"
$ws.Range("G34").Value = "10"
$ws.Range("G35").Value = "10"
$ws.Range("G36").Value = "10"
$ws.Range("G37").Value = "10"
$ws.Range("G38").Value = "10"
$ws.Range("D39").Value = "0.02170"
$ws.Range("E39").Value = "3.17%"
$ws.Range("G39").Value = "10"
$ws.Range("D40").Value = "0.05216"
$ws.Range("E40").Value = "1.32%"
$ws.Range("G40").Value = "10"
$ws.Range("D41").Value = "0.007402"
$ws.Range("E41").Value = "-0.87%"
$ws.Range("G41").Value = "10"
$ws.Range("D42").Value = "0.009775"
$ws.Range("E42").Value = "-3.56%"
$ws.Range("G42").Value = "10"
$ws.Range("D43").Value = "0.1370"
$ws.Range("E43").Value = "0.75%"
$ws.Range("G43").Value = "10"
$ws.Range("D44").Value = "0.002110"
$ws.Range("E44").Value = "-1.40%"
$ws.Range("G44").Value = "10"
$ws.Range("D45").Value = "0.009855"
$ws.Range("E45").Value = "7.33%"
$ws.Range("G45").Value = "10"
$ws.Range("D46").Value = "0.00006387"
$ws.Range("E46").Value = "1.73%"
$ws.Range("G46").Value = "10"
$ws.Range("D47").Value = "0.00000000747"
$ws.Range("E47").Value = "-0.35%"
$ws.Range("G47").Value = "10"
$ws.Range("D48").Value = "0.002758"
$ws.Range("E48").Value = "-10.04%"
$ws.Range("G48").Value = "10"
$ws.Range("D49").Value = "0.0009966"
$ws.Range("E49").Value = "-37.72%"
$ws.Range("G49").Value = "10"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("E50").Value = "-0.35%"
$ws.Range("G50").Value = "10"
$ws.Range("D51").Value = "0.0001993"
$ws.Range("E51").Value = "-0.35%"
$ws.Range("G51").Value = "10"
